# "add a minigame form" - insert a new row into the MainIcon table for the
# minigame panel button, ahead of the existing SideButton7 ("烹饪") row,
# shifting the following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row at row 23 (existing data rows 23-27 shift down to 24-28).
$ws.Rows("23:23").Insert()

# The table's own range doesn't auto-grow from a plain row insert, so resize
# it explicitly to cover the new row.
$lo.Resize($ws.Range("A1:L28"))

# Populate the new row with the minigame-panel entry.
$ws.Range("A23").Value = 42
$ws.Range("B23").Value = "游戏"
$ws.Range("C23").Value = "打开迷你游戏面板"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2
$ws.Range("L23").Value = "SideButton7"

# Match the author's final selection/cursor position.
$ws.Range("F23").Select()

# Recolor the workbook theme's light-1 (window background) swatch.
$wb.Theme.ThemeColorScheme.Item(2).RGB = 13494986
